$wb = $excel.ActiveWorkbook

# --- Update keymap values ---
# "lower" sheet: G5/H5 swap from mouse-click bindings previously on other
# sheets to the actual mouse click keycodes.
$wsLower = $wb.Worksheets.Item("lower")
$wsLower.Range("G5").Value = "mkp LCLK"
$wsLower.Range("H5").Value = "mkp RCLK"

# "raise" sheet: G5/H5 just get their stray trailing whitespace trimmed to
# plain "none".
$wsRaise = $wb.Worksheets.Item("raise")
$wsRaise.Range("G5").Value = "none"
$wsRaise.Range("H5").Value = "none"

# "adjust" sheet: G5/H5 pick up the RGB toggle / mute bindings that used to
# live on the "lower" sheet.
$wsAdjust = $wb.Worksheets.Item("adjust")
$wsAdjust.Range("G5").Value = "rgb_ug RGB_TOG             "
$wsAdjust.Range("H5").Value = "kp C_MUTE        "

# --- Update view state (selection / active sheet) to match ---
# Order matters: selecting a range activates its sheet, so the sheet that
# should end up active must be selected last.
$wb.Worksheets.Item("default").Range("H5").Select()
$wb.Worksheets.Item("lower").Range("H5").Select()
$wb.Worksheets.Item("adjust").Range("G5").Select()
$wb.Worksheets.Item("raise").Range("H4").Select()
